$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.986.74'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').Value = '2.220.86'
$ws.Range('E3').Value = '  -1.43%  '
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '243.03'
$ws.Range('E5').Value = '  -1.91%  '
$ws.Range('D6').Value = '0.628'
$ws.Range('E6').Value = '  +0.99%  '
$ws.Range('D7').Value = '74.29'
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.29%  '
$ws.Range('D10').Value = '44.23'
$ws.Range('E10').Value = '  +5.95%  '
$ws.Range('D11').Value = '0.0962'
$ws.Range('E11').Value = '  +2.41%  '
$ws.Range('D12').Value = '7.16'
$ws.Range('E12').Value = '  +0.83%  '
$ws.Range('E13').Value = '  +0.77%  '
$ws.Range('D14').Value = '2.552.71'
$ws.Range('E14').Value = '  -1.42%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '14.31'
$ws.Range('E15').Value = '  -1.41%  '
$ws.Range('B16').Value = 'Polygon'
$ws.Range('C16').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D16').Value = '0.847'
$ws.Range('E16').Value = '  -1.14%  '
$ws.Range('D17').Value = '2.205.63'
$ws.Range('E17').Value = '  -2.09%  '
$ws.Range('D18').Value = '41.904.82'
$ws.Range('E18').Value = '  -0.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000109'
$ws.Range('E19').Value = '  +11.54%  '
$ws.Range('D20').Value = '6.21'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('D21').Value = '72.47'
$ws.Range('E21').Value = '  +0.76%  '
$ws.Range('D22').Value = '11.07'
$ws.Range('E22').Value = '  +39.45%  '
$ws.Range('D23').Value = '229.96'
$ws.Range('E23').Value = '  -0.90%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.10'
$ws.Range('E24').Value = '  -7.75%  '
$ws.Range('B25').Value = 'Cosmos'
$ws.Range('C25').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.60'
$ws.Range('E25').Value = '  +3.45%  '
$ws.Range('B26').Value = 'Dai'
$ws.Range('C26').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.60'
$ws.Range('E27').Value = '  +1.05%  '
$ws.Range('E28').Value = '  -1.71%  '
$ws.Range('E29').Value = '  -3.35%  '
$ws.Range('D30').Value = '166.68'
$ws.Range('E30').Value = '  -1.49%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '20.60'
$ws.Range('E31').Value = '  -0.60%  '
$ws.Range('D32').Value = '5.68'
$ws.Range('E32').Value = '  +15.88%  '
$ws.Range('D33').Value = '0.0803'
$ws.Range('E33').Value = '  -2.59%  '
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('D35').Value = '29.52'
$ws.Range('E35').Value = '  -4.33%  '
$ws.Range('E36').Value = '  -4.32%  '
$ws.Range('E37').Value = '  -5.31%  '
$ws.Range('E38').Value = '  -1.80%  '
$ws.Range('D39').Value = '13.03'
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('D40').Value = '2.15'
$ws.Range('E40').Value = '  -1.65%  '
$ws.Range('D41').Value = '65.18'
$ws.Range('E41').Value = '  +5.60%  '
$ws.Range('E42').Value = '  -2.01%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.200'
$ws.Range('E43').Value = '  -1.50%  '
$ws.Range('D44').Value = '8.75'
$ws.Range('E44').Value = '  +0.89%  '
$ws.Range('D45').Value = '104.35'
$ws.Range('E45').Value = '  -3.40%  '
$ws.Range('E46').Value = '  +0.49%  '
$ws.Range('D47').Value = '2.45'
$ws.Range('E47').Value = '  +7.20%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.70'
$ws.Range('E50').Value = '  +0.36%  '
$ws.Range('D51').Value = '2.426.91'
$ws.Range('E51').Value = '  -1.46%  '
